# RELEASE: updated CHANGES.TXT and QR version
$d = $word.ActiveDocument

# 1. UVVM Utility Library version 2.14.0 -> 2.15.0
$found = $d.Content.Find.Execute("UVVM Utility Library (UVVM-Util), version 2.14.0 and up", $true, $false, $false, $false, $false, $true, 1, $false, "UVVM Utility Library (UVVM-Util), version 2.15.0 and up", 2)

# 2. UVVM VVC Framework version 2.10.0 -> 2.11.0
$found2 = $d.Content.Find.Execute("UVVM VVC Framework, version 2.10.0 and up", $true, $false, $false, $false, $false, $true, 1, $false, "UVVM VVC Framework, version 2.11.0 and up", 2)

# 3. Footer "Last update" date 2020-10-02 -> 2020-10-05
foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers.Item(1)
    if ($ftr.Exists) {
        $ftr.Range.Find.Execute("2020-10-02", $true, $false, $false, $false, $false, $true, 1, $false, "2020-10-05", 2)
    }
}
